# Fruta / hortaliza, semanal
# Insert 7 new daily rows (2022-03-08 / serial 44628) for Nectarín at
# Macroferia Regional de Talca, pushing the existing historical rows
# (625-635) down to (632-642).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 7 blank rows above the current row 625 ---
$ws.Rows("625:631").Insert()

# --- 2. Common (constant) column values shared by every row in this block ---
$colA = 5
$colB = "Macroferia Regional de Talca"
$colC = "Maule"
$colE = 7
$colF = "Fruta"
$colG = 100103
$colH = "Frutos de hueso (carozo)"
$colI = 100103006
$colJ = "Nectarín"
$colR = "Región de O'Higgins"
$newDate = 44628

# --- 3. Row-specific data for the 7 newly inserted rows (625-631) ---
$newRows = @(
    @{ Row=625; K="Artic Snow";     L="Especial";                M=160; N=14000; O=14000; P=14000; Q="$/bandeja 18 kilos granel";  S=778; T=18 },
    @{ Row=626; K="Artic Snow";     L="Extra (doble especial)";  M=180; N=16000; O=16000; P=16000; Q="$/bandeja 18 kilos granel";  S=889; T=18 },
    @{ Row=627; K="Artic Snow";     L="Primera";                 M=250; N=12000; O=12000; P=12000; Q="$/bandeja 18 kilos granel";  S=667; T=18 },
    @{ Row=628; K="Artic Snow";     L="Segunda";                 M=200; N=10000; O=10000; P=10000; Q="$/bandeja 18 kilos granel";  S=556; T=18 },
    @{ Row=629; K="June Pearl";     L="Extra (doble especial)";  M=150; N=16000; O=16000; P=16000; Q="$/bandeja 18 kilos granel";  S=889; T=18 },
    @{ Row=630; K="September Red";  L="Especial";                M=200; N=14000; O=14000; P=14000; Q="$/caja 15 kilos empedrada"; S=933; T=15 },
    @{ Row=631; K="Venus";          L="Extra (doble especial)";  M=120; N=16000; O=16000; P=16000; Q="$/bandeja 18 kilos granel";  S=889; T=18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $colA
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $newDate
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
    $ws.Cells.Item($row, 10).Value = $colJ
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $colR
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
